$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.182.00"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.58"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.62"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.26"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0593"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.781.17"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.562.08"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.78"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.84"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.150.42"
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.54"
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.13"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0463"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.385.73"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.949"
$ws.Range("E36").Value = "  -2.73%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.814"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.517"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.983"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.36"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.694.51"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.50"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  +0.11%  "
